$d = $word.ActiveDocument

# --- Change 1: remove the stray _GoBack bookmark after "1.2 - References" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2: "8- Programming Language " heading becomes "9- Programming Language ",
#     with a fresh _GoBack bookmark inserted right after the new "9" ---
$r = $d.Content
$found = $r.Find.Execute("8- Programming Language", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $r.Start
    # Place the bookmark right after the single digit, before the text is touched,
    # so the run split happens cleanly at that boundary.
    $bmRange = $d.Range($start + 1, $start + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    # Now swap the leading digit "8" -> "9".
    $digitRange = $d.Range($start, $start + 1)
    $digitRange.Text = "9"
}

# --- Change 3: header's cached PAGE field result "3" -> "6" ---
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)
